$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.457.29"
$ws.Range("E2").Value = "  -1.84%  "

$ws.Range("D3").Value = "1.849.05"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.99"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6414"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2988"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07450"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.22"
$ws.Range("E10").Value = "  -0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07630"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("D12").Value = "1.848.55"
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.015"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6830"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.47"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009492"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.134"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").Value = "29.492.82"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "2.068.46"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.82"
$ws.Range("E20").Value = "  -2.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -1.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.650"
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.94"
$ws.Range("E25").Value = "  -1.74%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1410"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.467"
$ws.Range("E27").Value = "  -1.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.74"
$ws.Range("E28").Value = "  -1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05973"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.255"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.120"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.049"
$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.868"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.173"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7182"
$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.600"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.794"
$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01774"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("D40").Value = "1.198.32"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9118"
$ws.Range("E41").Value = "  -2.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.162"
$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"

$ws.Range("D44").Value = "2.003.10"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.35"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.21"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.296"
$ws.Range("E47").Value = "  +8.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000121"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4024"
$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.051"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.656"
$ws.Range("E51").Value = "  +1.21%  "

Write-Host "Updated cryptos list"